# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.544.64'
$ws.Range("E2").Value = '  -0.90%  '
$ws.Range("D3").Value = '2.582.92'
$ws.Range("E3").Value = '  -1.84%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'" + '583.22'
$ws.Range("E5").Value = '  -1.50%  '
$ws.Range("D6").Value = "'" + '166.30'
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -1.60%  '
$ws.Range("D9").Value = '2.583.16'
$ws.Range("E9").Value = '  -1.83%  '
$ws.Range("D10").Value = "'" + '0.137'
$ws.Range("E10").Value = '  -3.91%  '
$ws.Range("E11").Value = '  +0.25%  '
$ws.Range("D12").Value = "'" + '0.356'
$ws.Range("E12").Value = '  -1.43%  '
$ws.Range("D13").Value = "'" + '5.17'
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("D14").Value = "'" + '26.70'
$ws.Range("E14").Value = '  -3.25%  '
$ws.Range("D15").Value = '3.053.48'
$ws.Range("E15").Value = '  -1.85%  '
$ws.Range("E16").Value = '  -2.21%  '
$ws.Range("D17").Value = '66.352.85'
$ws.Range("E17").Value = '  -1.00%  '
$ws.Range("D18").Value = '2.594.36'
$ws.Range("E18").Value = '  -0.98%  '
$ws.Range("D19").Value = "'" + '11.41'
$ws.Range("E19").Value = '  -5.80%  '
$ws.Range("D20").Value = "'" + '7.74'
$ws.Range("E20").Value = '  -3.64%  '
$ws.Range("D21").Value = "'" + '352.50'
$ws.Range("E21").Value = '  -1.98%  '
$ws.Range("D22").Value = "'" + '4.23'
$ws.Range("E22").Value = '  -2.74%  '
$ws.Range("D23").Value = "'" + '4.60'
$ws.Range("E23").Value = '  -1.55%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("E25").Value = '  -3.41%  '
$ws.Range("D26").Value = "'" + '10.03'
$ws.Range("E26").Value = '  -8.17%  '
$ws.Range("D27").Value = "'" + '68.88'
$ws.Range("E27").Value = '  -2.76%  '
$ws.Range("D28").Value = '2.714.98'
$ws.Range("E28").Value = '  -1.86%  '
$ws.Range("E29").Value = '  +0.33%  '
$ws.Range("D30").Value = '0.0₃0986'
$ws.Range("E30").Value = '  -2.53%  '
$ws.Range("D31").Value = "'" + '535.34'
$ws.Range("E31").Value = '  -3.37%  '
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("E33").Value = '  -2.68%  '
$ws.Range("E34").Value = '  -2.38%  '
$ws.Range("E35").Value = '  -1.66%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").Value = "'" + '1.46'
$ws.Range("E37").Value = '  -3.21%  '
$ws.Range("D38").Value = "'" + '156.68'
$ws.Range("E38").Value = '  -0.51%  '
$ws.Range("D39").Value = "'" + '18.75'
$ws.Range("E39").Value = '  -2.31%  '
$ws.Range("D40").Value = "'" + '0.360'
$ws.Range("E40").Value = '  -1.98%  '
$ws.Range("D41").Value = "'" + '18.24'
$ws.Range("E41").Value = '  +1.83%  '
$ws.Range("E42").Value = '  -0.53%  '
$ws.Range("D43").Value = "'" + '5.11'
$ws.Range("E43").Value = '  -1.99%  '
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").Value = "'" + '2.41'
$ws.Range("E45").Value = '  -2.35%  '
$ws.Range("D46").Value = '0.0₆0286'
$ws.Range("E46").Value = '  -4.79%  '
$ws.Range("D47").Value = "'" + '149.27'
$ws.Range("E47").Value = '  -1.82%  '
$ws.Range("D48").Value = "'" + '0.567'
$ws.Range("E48").Value = '  -3.37%  '
$ws.Range("D49").Value = "'" + '3.72'
$ws.Range("E49").Value = '  -2.47%  '
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("E51").Value = '  -1.68%  '
